$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 257.375
$ws.Range("I4").Value = 222.71428
$ws.Range("K4").Value = 222.71428
$ws.Range("M4").Value = -108.71428
$ws.Range("H9").Value = 101.6
$ws.Range("I9").Value = 101.6
$ws.Range("K9").Value = 101.6
$ws.Range("M9").Value = 67.40000000000001
$ws.Range("H10").Value = 3918
$ws.Range("J10").Value = 5625
$ws.Range("L10").Value = 5625
$ws.Range("N10").Value = -6211
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("K11").Value = 100
$ws.Range("M11").Value = 40
$ws.Range("H12").Value = 124
$ws.Range("I12").Value = 124.333336
$ws.Range("J12").Value = 123
$ws.Range("K12").Value = 124.333336
$ws.Range("L12").Value = 123
$ws.Range("M12").Value = 45.666664
$ws.Range("N12").Value = -463
$ws.Range("H33").Value = 210.05882
$ws.Range("I33").Value = 176.375
$ws.Range("K33").Value = 176.375
$ws.Range("M33").Value = 52.625
$ws.Range("H69").Value = 3500
$ws.Range("J69").Value = 3500
$ws.Range("L69").Value = 10500
$ws.Range("N69").Value = -12248
$ws.Range("H72").Value = 3500
$ws.Range("J72").Value = 3500
$ws.Range("L72").Value = 31500
$ws.Range("N72").Value = -40236

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1100
$ws.Range("J3").Value = 1100
$ws.Range("L3").Value = 1100
$ws.Range("N3").Value = -1330
$ws.Range("H32").Value = 8174.9165
$ws.Range("I32").Value = 8236.272000000001
$ws.Range("K32").Value = 8236.272000000001
$ws.Range("M32").Value = -7949.272000000001
$ws.Range("H39").Value = 1500
$ws.Range("I39").Value = 1500
$ws.Range("K39").Value = 1500
$ws.Range("M39").Value = -980
$ws.Range("H61").Value = 3651.75
$ws.Range("I61").Value = 3202.3333
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3202.3333
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2990.3333
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 2753
$ws.Range("I74").Value = 2004
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2004
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1130
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 2753
$ws.Range("I77").Value = 2004
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 10020
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -5652
$ws.Range("N77").Value = -33736
$ws.Range("H136").Value = 3651.75
$ws.Range("I136").Value = 3202.3333
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 9606.999899999999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -7056.999899999999
$ws.Range("N136").Value = -20100

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 747.6667
$ws.Range("I22").Value = 759
$ws.Range("J22").Value = 725
$ws.Range("K22").Value = 759
$ws.Range("L22").Value = 725
$ws.Range("M22").Value = -586
$ws.Range("N22").Value = -1071
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1550
$ws.Range("I31").Value = 1550
$ws.Range("K31").Value = 1550
$ws.Range("M31").Value = -1255
$ws.Range("H34").Value = 1550
$ws.Range("I34").Value = 1550
$ws.Range("K34").Value = 1550
$ws.Range("M34").Value = -1348
$ws.Range("H58").Value = 4803.6
$ws.Range("I58").Value = 2672.6667
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 2672.6667
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = -2469.6667
$ws.Range("N58").Value = -8406
$ws.Range("H105").Value = 3391.5557
$ws.Range("I105").Value = 2865.3333
$ws.Range("K105").Value = 2865.3333
$ws.Range("M105").Value = -1118.3333
$ws.Range("H132").Value = 3398.1
$ws.Range("I132").Value = 1330.3334
$ws.Range("K132").Value = 3991.0002
$ws.Range("M132").Value = -1461.0002
$ws.Range("H136").Value = 4803.6
$ws.Range("I136").Value = 2672.6667
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 8018.000100000001
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -5468.000100000001
$ws.Range("N136").Value = -29100

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20.555555
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 23
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 138
$ws.Range("M2").Value = 107
$ws.Range("N2").Value = -364
$ws.Range("H4").Value = 1008
$ws.Range("I4").Value = 600
$ws.Range("J4").Value = 1416
$ws.Range("K4").Value = 1800
$ws.Range("L4").Value = 4248
$ws.Range("M4").Value = -1688
$ws.Range("N4").Value = -4472
$ws.Range("H40").Value = 90.2
$ws.Range("I40").Value = 75
$ws.Range("J40").Value = 100.333336
$ws.Range("K40").Value = 300
$ws.Range("L40").Value = 401.333344
$ws.Range("M40").Value = -231
$ws.Range("N40").Value = -539.333344
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H134").Value = 500
$ws.Range("I134").Value = 500
$ws.Range("K134").Value = 1500
$ws.Range("M134").Value = 3570

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H80").Value = 1978.6
$ws.Range("I80").Value = 1332.6666
$ws.Range("J80").Value = 2947.5
$ws.Range("K80").Value = 1332.6666
$ws.Range("L80").Value = 2947.5
$ws.Range("M80").Value = -334.6666
$ws.Range("N80").Value = -4943.5
$ws.Range("H83").Value = 1978.6
$ws.Range("I83").Value = 1332.6666
$ws.Range("J83").Value = 2947.5
$ws.Range("K83").Value = 6663.333000000001
$ws.Range("L83").Value = 14737.5
$ws.Range("M83").Value = -1671.333000000001
$ws.Range("N83").Value = -24721.5
$ws.Range("H92").Value = 5881
$ws.Range("J92").Value = 5881
$ws.Range("L92").Value = 5881
$ws.Range("N92").Value = -9625
$ws.Range("H130").Value = 47498.75
$ws.Range("I130").Value = 19997.5
$ws.Range("J130").Value = 75000
$ws.Range("K130").Value = 19997.5
$ws.Range("L130").Value = 75000
$ws.Range("M130").Value = -14977.5
$ws.Range("N130").Value = -85040

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3803.25
$ws.Range("I61").Value = 3735.2
$ws.Range("J61").Value = 3916.6667
$ws.Range("K61").Value = 3735.2
$ws.Range("L61").Value = 3916.6667
$ws.Range("M61").Value = -3533.2
$ws.Range("N61").Value = -4320.6667
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 330.6
$ws.Range("I82").Value = 330.6
$ws.Range("K82").Value = 330.6
$ws.Range("M82").Value = 30.39999999999998
$ws.Range("H85").Value = 330.6
$ws.Range("I85").Value = 330.6
$ws.Range("K85").Value = 330.6
$ws.Range("M85").Value = 917.4
$ws.Range("H93").Value = 743.875
$ws.Range("I93").Value = 730.4
$ws.Range("J93").Value = 766.3333
$ws.Range("K93").Value = 730.4
$ws.Range("L93").Value = 766.3333
$ws.Range("M93").Value = 517.6
$ws.Range("N93").Value = -3262.3333
$ws.Range("H100").Value = 2775.75
$ws.Range("I100").Value = 2775.75
$ws.Range("K100").Value = 2775.75
$ws.Range("M100").Value = -2234.75
$ws.Range("H113").Value = 3803.25
$ws.Range("I113").Value = 3735.2
$ws.Range("J113").Value = 3916.6667
$ws.Range("K113").Value = 3735.2
$ws.Range("L113").Value = 3916.6667
$ws.Range("M113").Value = -1565.2
$ws.Range("N113").Value = -8256.6667

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1227.6666
$ws.Range("I126").Value = 1249.8572
$ws.Range("K126").Value = 3749.5716
$ws.Range("M126").Value = -1279.5716
